# Apply the latest MeritOrderLänder simplifications to the results sheet.
# Column A values are reduced by 1270, and column N values are set to a
# uniform 25400 (previously a uniform 24130) for data rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @{
    2  = 35154.0
    3  = 33264.0
    4  = 31573.0
    5  = 28188.0
    6  = 25983.0
    7  = 25465.0
    8  = 25687.0
    9  = 25993.0
    10 = 25667.0
    11 = 26257.0
}

foreach ($row in $colA.Keys) {
    $ws.Cells.Item($row, 1).Value = $colA[$row]
    $ws.Cells.Item($row, 14).Value = 25400.0
}
